$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '48.162.20'
$ws.Cells.Item(2, 5).Value = '  +1.89%  '

$ws.Cells.Item(3, 4).Value = '2.509.15'
$ws.Cells.Item(3, 5).Value = '  +0.68%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '321.07'
$ws.Cells.Item(5, 5).Value = '  -0.21%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '108.68'
$ws.Cells.Item(6, 5).Value = '  +0.29%  '

$ws.Cells.Item(7, 5).Value = '  +0.84%  '

$ws.Cells.Item(8, 5).Value = '  +0.00%  '

$ws.Cells.Item(9, 5).Value = '  +0.68%  '

$ws.Cells.Item(10, 5).Value = '  +1.83%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '20.20'
$ws.Cells.Item(11, 5).Value = '  +10.06%  '

$ws.Cells.Item(12, 5).Value = '  +1.05%  '

$ws.Cells.Item(13, 5).Value = '  +0.58%  '

$ws.Cells.Item(14, 5).Value = '  +0.95%  '

$ws.Cells.Item(15, 4).Value = '2.900.30'

$ws.Cells.Item(16, 4).Value = '2.510.19'
$ws.Cells.Item(16, 5).Value = '  +0.79%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.844'
$ws.Cells.Item(17, 5).Value = '  +0.08%  '

$ws.Cells.Item(18, 4).Value = '47.998.02'
$ws.Cells.Item(18, 5).Value = '  +1.74%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '13.14'
$ws.Cells.Item(19, 5).Value = '  -0.06%  '

$ws.Cells.Item(20, 5).Value = '  +0.42%  '

$ws.Cells.Item(21, 4).Value = '0.0₃0952'
$ws.Cells.Item(21, 5).Value = '  +1.55%  '

$ws.Cells.Item(22, 5).Value = '  +0.86%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '72.09'
$ws.Cells.Item(23, 5).Value = '  +2.41%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '275.84'
$ws.Cells.Item(24, 5).Value = '  +12.46%  '

$ws.Cells.Item(25, 5).Value = '  +0.36%  '

$ws.Cells.Item(26, 5).Value = '  -0.07%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '25.88'
$ws.Cells.Item(27, 5).Value = '  +0.65%  '

$ws.Cells.Item(28, 2).Value = 'Cosmos'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '9.84'
$ws.Cells.Item(28, 5).Value = '  -1.15%  '

$ws.Cells.Item(29, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '35.40'
$ws.Cells.Item(29, 5).Value = '  +2.35%  '

$ws.Cells.Item(30, 2).Value = 'Kaspa'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.137'
$ws.Cells.Item(30, 5).Value = '  -0.56%  '

$ws.Cells.Item(31, 2).Value = 'Toncoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '2.10'
$ws.Cells.Item(31, 5).Value = '  -7.54%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '49.23'
$ws.Cells.Item(32, 5).Value = '  -1.11%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '19.39'
$ws.Cells.Item(33, 5).Value = '  -4.00%  '

$ws.Cells.Item(34, 5).Value = '  +0.29%  '

$ws.Cells.Item(35, 5).Value = '  -0.02%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.0786'
$ws.Cells.Item(36, 5).Value = '  +0.28%  '

$ws.Cells.Item(37, 5).Value = '  -0.25%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '4.61'
$ws.Cells.Item(38, 5).Value = '  -3.13%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.95'
$ws.Cells.Item(39, 5).Value = '  +0.73%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '123.23'
$ws.Cells.Item(40, 5).Value = '  +4.14%  '

$ws.Cells.Item(41, 5).Value = '  +0.21%  '

$ws.Cells.Item(42, 5).Value = '  -0.40%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '21.57'
$ws.Cells.Item(43, 5).Value = '  -7.01%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.0305'
$ws.Cells.Item(44, 5).Value = '  +3.21%  '

$ws.Cells.Item(45, 4).Value = '1.997.76'
$ws.Cells.Item(45, 5).Value = '  -0.02%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '3.13'
$ws.Cells.Item(46, 5).Value = '  +3.31%  '

$ws.Cells.Item(47, 5).Value = '  +4.41%  '

$ws.Cells.Item(48, 5).Value = '  -0.81%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '9.02'
$ws.Cells.Item(49, 5).Value = '  -1.34%  '

$ws.Cells.Item(50, 5).Value = '  +2.09%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '79.77'
$ws.Cells.Item(51, 5).Value = '  +2.57%  '
